# "Generate Report for Handoff"
#
# This script updates the localization-status report after a handoff run:
#   - Rows 4-7 (the items that were previously "low" priority / not yet
#     handed off) now have Priority = "ht" on both the zh-cn and de-de
#     sheets, and their "Latest Handoff Datetime" (zh-cn) is refreshed.
#   - The "Latest HO Xliff Generate Date" shown on the Overview sheet for
#     those same rows (status = "Ready for handoff") is refreshed; this
#     same timestamp string also happens to back the de-de sheet's
#     "Latest Handoff Datetime" for rows 4-7, so it is updated there too.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# zh-cn: Priority low -> ht for rows 4-7
$wsZhCn.Range("E4:E7").Value = "ht"

# zh-cn: Latest Handoff Datetime refreshed for rows 4-7
$wsZhCn.Range("H4:H7").Value = "2016-08-17 18:30:44"

# de-de: Priority low -> ht for rows 4-7
$wsDeDe.Range("E4:E7").Value = "ht"

# de-de: Latest Handoff Datetime for rows 4-7 shares the same underlying
# timestamp as the Overview's "Latest HO Xliff Generate Date" below, so it
# is refreshed to the same new value.
$wsDeDe.Range("H4:H7").Value = "2016-08-17 18:30:50"

# Overview: Latest HO Xliff Generate Date refreshed for rows 4-7
$wsOverview.Range("G4:G7").Value = "2016-08-17 18:30:50"
